$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 522
$ws.Range("I15").Value = 522
$ws.Range("K15").Value = 1566
$ws.Range("M15").Value = -1397
$ws.Range("H51").Value = 35716.168
$ws.Range("J51").Value = 40059.8
$ws.Range("L51").Value = 40059.8
$ws.Range("N51").Value = -41027.8
$ws.Range("H94").Value = 649.25
$ws.Range("I94").Value = 649.25
$ws.Range("K94").Value = 649.25
$ws.Range("M94").Value = -198.25
$ws.Range("H100").Value = 4185
$ws.Range("I100").Value = 3359.625
$ws.Range("J100").Value = 4845.3
$ws.Range("K100").Value = 3359.625
$ws.Range("L100").Value = 4845.3
$ws.Range("M100").Value = -2818.625
$ws.Range("N100").Value = -5927.3
$ws.Range("H127").Value = 860.1429000000001
$ws.Range("I127").Value = 860.1429000000001
$ws.Range("K127").Value = 2580.4287
$ws.Range("M127").Value = 2379.5713
$ws.Range("H137").Value = 1565
$ws.Range("I137").Value = 1248.32
$ws.Range("J137").Value = 2696
$ws.Range("K137").Value = 3744.96
$ws.Range("L137").Value = 8088
$ws.Range("M137").Value = -1194.96
$ws.Range("N137").Value = -13188
$ws.Range("H140").Value = 114494
$ws.Range("J140").Value = 114494
$ws.Range("L140").Value = 114494
$ws.Range("N140").Value = -124854
$ws.Range("H141").Value = 3264.9443
$ws.Range("J141").Value = 6658
$ws.Range("L141").Value = 19974
$ws.Range("N141").Value = -30334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2050.75
$ws.Range("I97").Value = 2034.3334
$ws.Range("K97").Value = 2034.3334
$ws.Range("M97").Value = -1538.3334
$ws.Range("H138").Value = 88404.42999999999
$ws.Range("J138").Value = 88404.42999999999
$ws.Range("L138").Value = 88404.42999999999
$ws.Range("N138").Value = -98684.42999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I11").Value = 1066.6666
$ws.Range("J11").Value = 1005
$ws.Range("K11").Value = 1066.6666
$ws.Range("L11").Value = 1005
$ws.Range("M11").Value = -926.6666
$ws.Range("N11").Value = -1285
$ws.Range("H94").Value = 1100.5
$ws.Range("I94").Value = 1053.6875
$ws.Range("K94").Value = 1053.6875
$ws.Range("M94").Value = -602.6875
$ws.Range("H99").Value = 4274.25
$ws.Range("I99").Value = 4162.8184
$ws.Range("K99").Value = 4162.8184
$ws.Range("M99").Value = -2664.8184

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 59302.637
$ws.Range("J16").Value = 161299.75
$ws.Range("L16").Value = 161299.75
$ws.Range("N16").Value = -161873.75
$ws.Range("H22").Value = 839.4286
$ws.Range("J22").Value = 894.75
$ws.Range("L22").Value = 894.75
$ws.Range("N22").Value = -1594.75
$ws.Range("H62").Value = 2944.8333
$ws.Range("I62").Value = 2920.25
$ws.Range("K62").Value = 2920.25
$ws.Range("M62").Value = -2296.25
$ws.Range("H65").Value = 2944.8333
$ws.Range("I65").Value = 2920.25
$ws.Range("K65").Value = 14601.25
$ws.Range("M65").Value = -11481.25
$ws.Range("H86").Value = 30310726
$ws.Range("I86").Value = 66672740
$ws.Range("J86").Value = 9049.833000000001
$ws.Range("K86").Value = 66672740
$ws.Range("L86").Value = 9049.833000000001
$ws.Range("M86").Value = -66671617
$ws.Range("N86").Value = -11295.833
$ws.Range("H89").Value = 30310726
$ws.Range("I89").Value = 66672740
$ws.Range("J89").Value = 9049.833000000001
$ws.Range("K89").Value = 333363700
$ws.Range("L89").Value = 45249.165
$ws.Range("M89").Value = -333358084
$ws.Range("N89").Value = -56481.165
$ws.Range("H113").Value = 59302.637
$ws.Range("J113").Value = 161299.75
$ws.Range("L113").Value = 161299.75
$ws.Range("N113").Value = -165639.75
$ws.Range("H122").Value = 133995.1
$ws.Range("I122").Value = 146266.05
$ws.Range("K122").Value = 438798.15
$ws.Range("M122").Value = -436348.15
$ws.Range("H130").Value = 63176.5
$ws.Range("J130").Value = 69332.336
$ws.Range("L130").Value = 69332.336
$ws.Range("N130").Value = -79372.336
$ws.Range("H132").Value = 6699.2163
$ws.Range("I132").Value = 5421.1304
$ws.Range("J132").Value = 8798.929
$ws.Range("K132").Value = 16263.3912
$ws.Range("L132").Value = 26396.787
$ws.Range("M132").Value = -13733.3912
$ws.Range("N132").Value = -31456.787
$ws.Range("H134").Value = 6395.1577
$ws.Range("I134").Value = 6157
$ws.Range("K134").Value = 18471
$ws.Range("M134").Value = -15936
$ws.Range("H141").Value = 74427.14
$ws.Range("J141").Value = 71832.164
$ws.Range("L141").Value = 71832.164
$ws.Range("N141").Value = -82192.164

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2557.4167
$ws.Range("I3").Value = 2146.913
$ws.Range("J3").Value = 11999
$ws.Range("K3").Value = 6440.739
$ws.Range("L3").Value = 35997
$ws.Range("M3").Value = -6328.739
$ws.Range("N3").Value = -36221
$ws.Range("H12").Value = 326.79486
$ws.Range("I12").Value = 383.45456
$ws.Range("J12").Value = 304.5357
$ws.Range("K12").Value = 1150.36368
$ws.Range("L12").Value = 913.6071000000001
$ws.Range("M12").Value = -977.3636799999999
$ws.Range("N12").Value = -1259.6071
$ws.Range("H133").Value = 5299.5
$ws.Range("I133").Value = 5299.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 15898.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -10838.5
$ws.Range("N133").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 50000
$ws.Range("J34").Value = 50000
$ws.Range("L34").Value = 50000
$ws.Range("N34").Value = -50536
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50630
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52184
$ws.Range("H80").Value = 4949.5
$ws.Range("I80").Value = 4900
$ws.Range("J80").Value = 4999
$ws.Range("K80").Value = 4900
$ws.Range("L80").Value = 4999
$ws.Range("M80").Value = -3902
$ws.Range("N80").Value = -6995
$ws.Range("H83").Value = 4949.5
$ws.Range("I83").Value = 4900
$ws.Range("J83").Value = 4999
$ws.Range("K83").Value = 24500
$ws.Range("L83").Value = 24995
$ws.Range("M83").Value = -19508
$ws.Range("N83").Value = -34979
$ws.Range("H102").Value = 2699.4
$ws.Range("I102").Value = 2699.4
$ws.Range("K102").Value = 2699.4
$ws.Range("M102").Value = -1077.4
$ws.Range("H122").Value = 3596.6667
$ws.Range("I122").Value = 3348.5
$ws.Range("J122").Value = 4093
$ws.Range("K122").Value = 10045.5
$ws.Range("L122").Value = 12279
$ws.Range("M122").Value = -7595.5
$ws.Range("N122").Value = -17179
$ws.Range("H132").Value = 6961.913
$ws.Range("I132").Value = 6566.294
$ws.Range("K132").Value = 19698.882
$ws.Range("M132").Value = -17168.882
$ws.Range("H141").Value = 63798.4
$ws.Range("J141").Value = 63798.4
$ws.Range("L141").Value = 63798.4
$ws.Range("N141").Value = -74158.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16698.111
$ws.Range("I7").Value = 18330.143
$ws.Range("J7").Value = 10986
$ws.Range("K7").Value = 18330.143
$ws.Range("L7").Value = 10986
$ws.Range("M7").Value = -18218.143
$ws.Range("N7").Value = -11210
$ws.Range("H40").Value = 25566.25
$ws.Range("I40").Value = 25029.6
$ws.Range("J40").Value = 28249.5
$ws.Range("K40").Value = 25029.6
$ws.Range("L40").Value = 28249.5
$ws.Range("M40").Value = -24893.6
$ws.Range("N40").Value = -28521.5
$ws.Range("H68").Value = 5991.375
$ws.Range("I68").Value = 4663
$ws.Range("J68").Value = 6788.4
$ws.Range("K68").Value = 4663
$ws.Range("L68").Value = 6788.4
$ws.Range("M68").Value = -3914
$ws.Range("N68").Value = -8286.4
$ws.Range("H71").Value = 5991.375
$ws.Range("I71").Value = 4663
$ws.Range("J71").Value = 6788.4
$ws.Range("K71").Value = 23315
$ws.Range("L71").Value = 33942
$ws.Range("M71").Value = -19571
$ws.Range("N71").Value = -41430
$ws.Range("H93").Value = 1739.375
$ws.Range("I93").Value = 1739.375
$ws.Range("K93").Value = 1739.375
$ws.Range("M93").Value = -491.375
$ws.Range("H126").Value = 16698.111
$ws.Range("I126").Value = 18330.143
$ws.Range("J126").Value = 10986
$ws.Range("K126").Value = 54990.429
$ws.Range("L126").Value = 32958
$ws.Range("M126").Value = -52520.429
$ws.Range("N126").Value = -37898
$ws.Range("H128").Value = 84999
$ws.Range("J128").Value = 84999
$ws.Range("L128").Value = 84999
$ws.Range("N128").Value = -94959
$ws.Range("H132").Value = 6210.6
$ws.Range("I132").Value = 5221.6
$ws.Range("J132").Value = 7199.6
$ws.Range("K132").Value = 15664.8
$ws.Range("L132").Value = 21598.8
$ws.Range("M132").Value = -13134.8
$ws.Range("N132").Value = -26658.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 11999
$ws.Range("J31").Value = 11999
$ws.Range("L31").Value = 11999
$ws.Range("N31").Value = -12695
$ws.Range("H122").Value = 2903.0833
$ws.Range("J122").Value = 2628.6667
$ws.Range("L122").Value = 7886.000100000001
$ws.Range("N122").Value = -12786.0001
$ws.Range("H135").Value = 90219
$ws.Range("J135").Value = 90219
$ws.Range("L135").Value = 90219
$ws.Range("N135").Value = -100359
$ws.Range("H140").Value = 51473.332
$ws.Range("J140").Value = 52711
$ws.Range("L140").Value = 52711
$ws.Range("N140").Value = -63071
$ws.Range("H141").Value = 65666.336
$ws.Range("J141").Value = 65666.336
$ws.Range("L141").Value = 65666.336
$ws.Range("N141").Value = -76026.336
